$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new product row (row 3)
$ws.Range("A3").Value = "SKU002"
$ws.Range("B3").Value = "Apple Juice"

# Date cells - copy formatting (and number format style) from the row above
# so they reuse the same cell style (short date) instead of creating a new one.
$ws.Range("C2:E2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$ws.Range("C3").Value = 45987
$ws.Range("D3").Value = 45988
$ws.Range("E3").Value = 46157

$ws.Range("F3").Value = '"Milk; Sugar; Banana"'
$ws.Range("G3").Value = 9.5
$ws.Range("H3").Value = 8901234567891
$ws.Range("I3").Value = "D:\Routisync\logo.png"

# Update the LogoFile value for the existing Orange Juice row (row 2):
# "D:\Routisync\rsynclogo.png" -> "D:\Routisync\logo.png"
$ws.Range("I2").Value = "D:\Routisync\logo.png"

# Widen column H to fit the new (longer) logo path text
# (target stored width is 34.7265625; the engine quantizes ColumnWidth to
# pixel-based increments, so 33.8 is the input that lands closest to it)
$ws.Columns.Item(8).ColumnWidth = 33.8

# Update the active selection shown when the workbook is reopened
$ws.Range("H9").Select() | Out-Null
